$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Window geometry (best-effort - provide a larger canvas / narrower window)
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 500
$win.ScrollRow = 11
$win.ScrollColumn = 1

# ---------------------------------------------------------------------------
# 2. Sensor co-ordinate table updates (rows 9-15, columns G & H)
#    Several cells lose their shared formula and become plain literal
#    values, while others gain/keep a recalculated formula.
# ---------------------------------------------------------------------------

# l_ear (row 9) - both G & H become static values
$ws.Range("G9").Value = 1158
$ws.Range("H9").Value = 890

# r_ear (row 10) - both G & H become static values
$ws.Range("G10").Value = 1242
$ws.Range("H10").Value = 890

# left (row 11) - G becomes static, H keeps/gains formula
$ws.Range("G11").Value = 1120
$ws.Range("H11").Formula = "=1200-C11"

# bl_corner (row 12) - G keeps/gains formula, H becomes static
$ws.Range("G12").Formula = "=1200-D12"
$ws.Range("H12").Value = 1411

# tail (row 13) - G keeps/gains formula, H becomes static
$ws.Range("G13").Formula = "=1200-D13"
$ws.Range("H13").Value = 1430

# br_corner (row 14) - G keeps/gains formula, H becomes static
$ws.Range("G14").Formula = "=1200-D14"
$ws.Range("H14").Value = 1411

# right (row 15) - G becomes static, H keeps/gains formula
$ws.Range("G15").Value = 1280
$ws.Range("H15").Formula = "=1200-C15"

# ---------------------------------------------------------------------------
# 3. New cell F19 holding the raw sensor JSON payload (as text, with a
#    literal leading apostrophe character that must be preserved, and not
#    treated by Excel as a "force text" quote-prefix marker).
#    Building it as a formula and then converting the formula to a static
#    value via copy / paste-values keeps the apostrophe as real text
#    without introducing a quotePrefix cell style.
# ---------------------------------------------------------------------------
$jsonPayload = '[{"sensorName":"l_ear","x": 1158,"y":890,"angle":999},{"sensorName":"r_ear","x": 1242,"y":890,"angle":999},{"sensorName":"left","x": 1120,"y":1343,"angle":180},{"sensorName":"bl_corner","x": 1152,"y":1411,"angle":225},{"sensorName":"tail","x": 1200,"y":1430,"angle":270},{"sensorName":"br_corner","x": 1248,"y":1411,"angle":315},{"sensorName":"right","x": 1280,"y":1343,"angle":0}]'
$escapedPayload = $jsonPayload -replace '"', '""'
$payloadFormula = '="' + "'" + $escapedPayload + '"'

$ws.Range("F19").Formula = $payloadFormula
$ws.Range("F19").Copy()
$ws.Range("F19").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. Selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("G23").Select()
